{"js": "// The progress-report table has a trailing block of empty rows after the\n// last filled-in data row (\"- Vi\u1ebft b\u00e1o c\u00e1o\"). Trim the table down by\n// removing five of those empty rows, leaving the rest of the table intact.\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Load the text of every row so we can locate the \"- Vi\u1ebft b\u00e1o c\u00e1o\" row and\n// the run of empty rows that immediately follows it.\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < rows.items.length; i++) {\n  rows.items[i].load(\"values\");\n}\nawait context.sync();\n\nlet anchorIndex = -1;\nfor (let i = 0; i < rows.items.length; i++) {\n  const rowValues = rows.items[i].values[0] || [];\n  if (rowValues.some((cell) => (cell || \"\").indexOf(\"Vi\u1ebft b\u00e1o c\u00e1o\") !== -1)) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the '- Vi\u1ebft b\u00e1o c\u00e1o' row.\");\n}\n\n// Remove the 5 empty rows immediately following the anchor row (keep the\n// remaining empty rows at the end of the table untouched). Re-fetch the\n// row collection and sync after every single delete so row indices never\n// go stale while the deletes are applied.\nconst rowsToRemove = 5;\nconst firstEmptyIndex = anchorIndex + 1;\n\nfor (let i = 0; i < rowsToRemove; i++) {\n  const freshRows = table.rows;\n  freshRows.load(\"items\");\n  await context.sync();\n\n  freshRows.items[firstEmptyIndex].delete();\n  await context.sync();\n}\n", "ps1": "# The progress-report table ends with a long run of empty rows after the\n# last filled-in data row (\"- Vi\u1ebft b\u00e1o c\u00e1o\"). Trim the table by removing\n# five of those empty rows, leaving the remaining empty rows untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Locate the row that contains \"- Vi\u1ebft b\u00e1o c\u00e1o\" in its second column.\n$anchorRow = -1\nfor ($i = 1; $i -le $t.Rows.Count; $i++) {\n    $row = $t.Rows.Item($i)\n    $txt = $row.Cells.Item(2).Range.Text\n    if ($txt -like \"*Vi\u1ebft b\u00e1o c\u00e1o*\") {\n        $anchorRow = $i\n        break\n    }\n}\n\nif ($anchorRow -eq -1) {\n    throw \"Could not find the '- Vi\u1ebft b\u00e1o c\u00e1o' row.\"\n}\n\n# Delete the 5 empty rows immediately following the anchor row. Deleting\n# the same index repeatedly shifts the following rows up into that slot.\n$rowsToRemove = 5\nfor ($n = 0; $n -lt $rowsToRemove; $n++) {\n    $t.Rows.Item($anchorRow + 1).Delete()\n}\n"}
